$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "90.503.47"
$ws.Range("E2").Value = "  -0.40%  "

# Row 3
$ws.Range("D3").Value = "3.112.25"
$ws.Range("E3").Value = "  -1.70%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.25%  "

# Row 5
$ws.Range("D5").Value = "234.02"
$ws.Range("E5").Value = "  +8.45%  "

# Row 6
$ws.Range("D6").Value = "623.35"
$ws.Range("E6").Value = "  -0.72%  "

# Row 7
$ws.Range("D7").Value = "1.07"
$ws.Range("E7").Value = "  -5.13%  "

# Row 8
$ws.Range("D8").Value = "0.367"
$ws.Range("E8").Value = "  -1.91%  "

# Row 9
$ws.Range("E9").Value = "  +0.09%  "

# Row 10
$ws.Range("D10").Value = "3.107.02"
$ws.Range("E10").Value = "  -1.70%  "

# Row 11
$ws.Range("D11").Value = "0.737"
$ws.Range("E11").Value = "  -2.08%  "

# Row 12
$ws.Range("E12").Value = "  -2.89%  "

# Row 13
$ws.Range("D13").Value = "0.0000252"
$ws.Range("E13").Value = "  +2.32%  "

# Row 14
$ws.Range("D14").Value = "35.97"
$ws.Range("E14").Value = "  +2.73%  "

# Row 15
$ws.Range("D15").Value = "5.48"
$ws.Range("E15").Value = "  -3.67%  "

# Row 16
$ws.Range("D16").Value = "90.245.50"
$ws.Range("E16").Value = "  -0.34%  "

# Row 17
$ws.Range("E17").Value = "  -2.11%  "

# Row 18
$ws.Range("B18").Value = "SuiNetwork"
$ws.Range("C18").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D18").Value = "3.89"
$ws.Range("E18").Value = "  +5.36%  "

# Row 19
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.095.09"
$ws.Range("E19").Value = "  -1.47%  "

# Row 20
$ws.Range("D20").Value = "0.0000222"
$ws.Range("E20").Value = "  +4.66%  "

# Row 21
$ws.Range("D21").Value = "14.04"
$ws.Range("E21").Value = "  -2.18%  "

# Row 22
$ws.Range("D22").Value = "437.58"
$ws.Range("E22").Value = "  -6.30%  "

# Row 23
$ws.Range("D23").Value = "5.57"
$ws.Range("E23").Value = "  +5.62%  "

# Row 24
$ws.Range("D24").Value = "8.95"
$ws.Range("E24").Value = "  -1.84%  "

# Row 25
$ws.Range("E25").Value = "  +3.16%  "

# Row 26
$ws.Range("D26").Value = "7.58"
$ws.Range("E26").Value = "  -1.44%  "

# Row 27
$ws.Range("D27").Value = "88.75"
$ws.Range("E27").Value = "  -5.32%  "

# Row 28
$ws.Range("D28").Value = "12.09"
$ws.Range("E28").Value = "  -1.06%  "

# Row 29
$ws.Range("D29").Value = "3.245.27"
$ws.Range("E29").Value = "  -2.17%  "

# Row 30
$ws.Range("E30").Value = "  -0.04%  "

# Row 31
$ws.Range("E31").Value = "  +1.70%  "

# Row 32
$ws.Range("D32").Value = "0.160"
$ws.Range("E32").Value = "  -1.07%  "

# Row 33
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -2.67%  "

# Row 34
$ws.Range("D34").Value = "0.195"

# Row 35
$ws.Range("D35").Value = "25.87"
$ws.Range("E35").Value = "  -4.80%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.153"
$ws.Range("E36").Value = "  +7.79%  "

# Row 37
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "3.83"
$ws.Range("E37").Value = "  +4.60%  "

# Row 38
$ws.Range("E38").Value = "  +5.21%  "

# Row 39
$ws.Range("D39").Value = "504.06"
$ws.Range("E39").Value = "  -3.71%  "

# Row 40
$ws.Range("D40").Value = "1.90"
$ws.Range("E40").Value = "  -1.40%  "

# Row 41
$ws.Range("E41").Value = "  -1.87%  "

# Row 42
$ws.Range("D42").Value = "0.0898"
$ws.Range("E42").Value = "  +4.79%  "

# Row 43
$ws.Range("E43").Value = "  -0.15%  "

# Row 44
$ws.Range("D44").Value = "0.406"
$ws.Range("E44").Value = "  -2.29%  "

# Row 45
$ws.Range("E45").Value = "  +0.02%  "

# Row 46
$ws.Range("D46").Value = "3.45"
$ws.Range("E46").Value = "  +55.47%  "

# Row 47
$ws.Range("D47").Value = "1.90"
$ws.Range("E47").Value = "  -3.66%  "

# Row 48
$ws.Range("D48").Value = "0.692"
$ws.Range("E48").Value = "  +1.53%  "

# Row 49
$ws.Range("D49").Value = "152.60"
$ws.Range("E49").Value = "  +1.60%  "

# Row 50
$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").Value = "1.34"
$ws.Range("E50").Value = "  -0.99%  "

# Row 51
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").Value = "44.59"
$ws.Range("E51").Value = "  -1.70%  "
